$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.10982014941406522
$ws.Range("D2").Value = 3603.181185708

$ws.Range("C3").Value = 0.1410716416367616
$ws.Range("D3").Value = 3600.319577146

$ws.Range("B4").Value = -31.005961887755465
$ws.Range("C4").Value = 0.118061786714723
$ws.Range("D4").Value = 3600.228164379

$ws.Range("C5").Value = 0.10316706501824413
$ws.Range("D5").Value = 3600.3665329

$ws.Range("D6").Value = 219.9979068

$ws.Range("C7").Value = 0.1302676711416746
$ws.Range("D7").Value = 3600.399690093

$ws.Range("C8").Value = 0.17110576102726616
$ws.Range("D8").Value = 3600.32528273

$ws.Range("C9").Value = 0.10588533590270065
$ws.Range("D9").Value = 3600.313373911

$ws.Range("C10").Value = 0.14952658509332356
$ws.Range("D10").Value = 3600.314417251

$ws.Range("C11").Value = 0.1608880825972639
$ws.Range("D11").Value = 3600.319067856
